$d = $word.ActiveDocument

# The document contains several "<id>...</id>" markers, each originally
# split across three differently-formatted runs:
#   run1 "<id>"   (Courier New, color 7f6000, sz 18)
#   run2 "<value>" (Arial, color 000000, sz 22)
#   run3 "</id>"  (Courier New, color 7f6000, sz 18)
# Two of them ("p033r_1" and "p033r_2" - the plain, non "fig_" ids) need
# to be collapsed into a single run (keeping run1's formatting) whose
# text is the concatenation "<id>value</id>".
#
# We locate each target paragraph by its exact rendered text, then
# rebuild it: insert the full combined text right after the opening
# "<id>" fragment (which keeps that run's formatting/attributes, incl.
# xml:space="preserve"), and delete the remainder of the paragraph's
# text that followed it.

function Merge-IdRun([string]$innerValue) {
    $target = "<id>" + $innerValue + "</id>"
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq $target -or $p.Range.Text -eq ($target + [char]13)) {
            $full = $p.Range
            $start = $full.Start
            $openLen = 4   # length of "<id>"
            $r1 = $d.Range($start, $start + $openLen)
            $rRest = $d.Range($start + $openLen, $start + $target.Length)
            $rRest.Text = ""
            $r1.InsertAfter($innerValue + "</id>")
            return
        }
    }
}

Merge-IdRun "p033r_1"
Merge-IdRun "p033r_2"
